# Bank Account Change Letter config.xlsx update
# - wording updates for LogFilePath and WorkpackageName values (breakpoints/title wording)
# - remove no-longer-needed config rows (CaseTitle, CaseOrigin, CaseType, CorrectFormType, LetterAttachedTitle)
# - table/dimension auto-shrinks from A1:C35 to A1:C30
# - scroll/selection moved to B18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wording of the log file path and workpackage name
$ws.Range("B18").Value = "C:\Users\{0}\Desktop\Bank Change Letter Logs_{1}.xlsx"
$ws.Range("B19").Value = "RPA Bank Change Letter"

# Remove rows that are no longer needed (delete bottom-up so row numbers stay valid)
# Row 30: LetterAttachedTitle
$ws.Rows(30).Delete()
# Row 27: CorrectFormType
$ws.Rows(27).Delete()
# Row 26: CaseType
$ws.Rows(26).Delete()
# Row 25: CaseOrigin
$ws.Rows(25).Delete()
# Row 24: CaseTitle
$ws.Rows(24).Delete()

# Update the saved view/selection to match where editing finished
$ws.Application.Goto($ws.Range("A18"), $true)
$ws.Range("B18").Select()

$wb.Save()
